# Applies the content changes described in the diff:
#  - TC1 (rows 10-12), TC2 (rows 20-22), TC3 (rows 30-32) blocks each gain
#    an extra "fill fields / submit" step + failure-message step before the
#    final "pick suggested user / submit" success step, and the opening
#    step's wording changes from "inicia" to "abre".
#
# Concretely (per test-case block, using B = Steps column, D = Expected
# Results column):
#   Step 1 (row x0): B text changes "inicia" -> "abre" (D unchanged)
#   Step 2 (row x1): B becomes "Usuario do Sistema preenche os campos e
#                     clica no botao entrar"; D becomes the TC-specific
#                     failure message that used to live in the old step 2
#   Step 3 (row x2): B becomes "Usuario do Sistema seleciona um nome de
#                     usuario sugerido, digita a senha e clica no botao
#                     entrar" (the old step 3 text); D stays "SYSTEM exibe
#                     uma mensagem de sucesso"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$blocks = @(
    @{ Step1 = 10; Step2 = 11; Step3 = 12; Fail = "SYSTEM alerta que o nome de usuario e/ou senha estao incorretos" },
    @{ Step1 = 20; Step2 = 21; Step3 = 22; Fail = "SYSTEM alerta que o TJSeg (sistema que fornece as permissoes de acesso e escrita) esta fora do ar" },
    @{ Step1 = 30; Step2 = 31; Step3 = 32; Fail = "SYSTEM alerta que o CAS (sistema de autorizacao login-senha) esta fora do ar" }
)

foreach ($block in $blocks) {
    $r1 = $block.Step1
    $r2 = $block.Step2
    $r3 = $block.Step3

    # Step 1: wording tweak only
    $ws.Range("B$r1").Value = "Usuario do Sistema abre a tela de login atraves da opcao de Login no canto superior direito"

    # Step 2: becomes the "fill fields and submit" step, expecting the
    # test-case-specific failure message
    $ws.Range("B$r2").Value = "Usuario do Sistema preenche os campos e clica no botao entrar"
    $ws.Range("D$r2").Value = $block.Fail

    # Step 3: becomes the "pick suggested user, submit" success step
    $ws.Range("B$r3").Value = "Usuario do Sistema seleciona um nome de usuario sugerido, digita a senha e clica no botao entrar"
}
